$d = $word.ActiveDocument

# The original "_GoBack" bookmark sits right after the last character of
# what used to be the final verse line ("...corredor."). It needs to move
# to mark the new last-edited spot: right after the new title text. Drop
# it now; it gets re-created below, in its new home.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a new paragraph at the very start of the document that becomes
# the poem's title. A trailing placeholder character ("X") is included
# temporarily -- it keeps the bookmark position we add further down from
# landing exactly on the paragraph's end boundary (which this Word engine
# mishandles for zero-length bookmarks), and is removed right after.
$ins = $d.Range(0, 0)
$ins.InsertBefore("TÍTULOX`r")

# The freshly inserted text is now the document's first paragraph; restyle
# it as the title: distinctive "Heavitas" font (ascii/hAnsi only -- the cs
# / complex-script font stays "Open Sans"), same purple color, 24pt text
# (sz/szCs = 48 half-points).
$titlePara = $d.Paragraphs(1)
$titleRun = $titlePara.Range
$titleRun.Font.Name = "Heavitas"
$titleRun.Font.Size = 24
$titleRun.Font.SizeBi = 24

# Re-create "_GoBack" right after "TÍTULO" (before the placeholder "X").
$bmRange = $d.Range(6, 6)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholder = $d.Range(6, 7)
$placeholder.Delete()
